# finish PHP1 notif improvement (case_study_app_V3)

$wb = $excel.ActiveWorkbook

$wsToDo = $wb.Worksheets.Item("To-Do List")
$wsCase = $wb.Worksheets.Item("todo case studies app")
$wsTreat = $wb.Worksheets.Item("todo treatments app")

# --- Sheet "To-Do List" : table ToDoList (rows 4-7) ---
$wsToDo.Range("G4").Value = 0.75
$wsToDo.Range("G6").Value = 0.5

# --- Sheet "To-Do List" : table ToDoList3 (rows 11-14) ---
$wsToDo.Range("D11").Value = "Complete"
$wsToDo.Range("G11").Value = 1
$wsToDo.Range("D13").Value = "Deferred"
$wsToDo.Range("I13").Value = "Cancelled"
$wsToDo.Range("D14").Value = "Complete"
$wsToDo.Range("G14").Value = 1
$wsToDo.Range("I14").Value = "Done"
$wsToDo.Range("I12").Value = "Not necessary at this time"

# --- Sheet "todo case studies app" : table ToDoList32 (rows 3-6) ---
$wsCase.Range("C3").Value = "In Progress"
$wsCase.Range("F3").Value = 0.5
$wsCase.Range("C4").Value = "Deferred"
$wsCase.Range("F4").ClearContents()
$wsCase.Range("H4").Value = "Canceled"

# --- Back to Sheet "To-Do List" for the last brand new string ---
$wsToDo.Range("I7").Value = "When the design is frozen"

# --- Sheet "todo treatments app" : table ToDoList324 (rows 3-4) ---
$wsTreat.Range("B3").Value = "Low"
$wsTreat.Range("C3").Value = "Not Started"
$wsTreat.Range("B4").Value = "Low"
$wsTreat.Range("C4").Value = "Not Started"

# Move active selections to match final workbook state
$wsToDo.Range("L6").Select()
$wsCase.Range("D12").Select()
$wsTreat.Range("B4").Select()

# Keep "To-Do List" as the active/front sheet (unchanged by this edit)
$wsToDo.Select()
